$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml($paraIndex, $innerXml, $isLast) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $startPos = $full.Start
    [void]$full.Delete()
    $ins = $d.Range($startPos, $startPos)
    if ($isLast) {
        $xml = "<w:p $wns>$innerXml</w:p>"
        [void]$ins.InsertXML($xml)
    } else {
        # Insert the real paragraph plus a throwaway empty paragraph so a
        # paragraph break is actually created, then merge the throwaway
        # paragraph mark forward to restore the original paragraph count.
        $xml = "<w:p $wns>$innerXml</w:p><w:p $wns></w:p>"
        [void]$ins.InsertXML($xml)
        $emptyP = $d.Paragraphs($paraIndex + 1)
        $mark = $d.Range($emptyP.Range.Start, $emptyP.Range.End)
        [void]$mark.Delete()
    }
}

$customerInner = '<w:r><w:t>customer (</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>username</w:t></w:r><w:r><w:t xml:space="preserve">, pin, </w:t></w:r><w:r><w:t>f</w:t></w:r><w:r><w:t xml:space="preserve">name, </w:t></w:r><w:r><w:t>l</w:t></w:r><w:r><w:t xml:space="preserve">name, street, city, state, zip, </w:t></w:r><w:r><w:t>card_</w:t></w:r><w:r><w:t xml:space="preserve">type, </w:t></w:r><w:r><w:t>card_num</w:t></w:r><w:r><w:t>, exp_date)</w:t></w:r>'

$bookInner = '<w:r><w:t>book (</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>isbn</w:t></w:r><w:r><w:t>, title, author, publisher, price, genre)</w:t></w:r>'

$purchaseInner = '<w:r><w:t>purchase (</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>purchaseID</w:t></w:r><w:r><w:t>, total)</w:t></w:r>'

$ordersInner = '<w:r><w:t>order</w:t></w:r><w:r><w:t>s</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:i/><w:u w:val="single"/></w:rPr><w:t>username, isbn, purchaseID</w:t></w:r><w:r><w:t>, quantity)</w:t></w:r>'

$reviewInner = '<w:r><w:t>review (</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>reviewID</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>isbn</w:t></w:r><w:r><w:t>, content)</w:t></w:r>'

Replace-ParagraphXml 6 $customerInner $false
Replace-ParagraphXml 7 $bookInner $false
Replace-ParagraphXml 8 $purchaseInner $false
Replace-ParagraphXml 9 $ordersInner $false
Replace-ParagraphXml 10 $reviewInner $true
